# Apply the activity-log update described in the commit:
# "Compiled successfully on Quartus and ModelSim. Updated logs"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log - Part 3")

# Fix a typo in the existing log entry on row 12 (drop trailing period, fix spelling)
$ws.Range("G12").Value = "Started working on the 64-barrel shifters but quickly got stuck on understanding on implementation"

# Fill in the four new activity-log entries on rows 28-31
$ws.Range("B28").Value = 6977
$ws.Range("C28").Value = 43933
$ws.Range("D28").Value = 0.84444444444444444
$ws.Range("E28").Value = 0.84583333333333333
$ws.Range("G28").Value = "Compiled all .vhd files on ModelSim and Quartus."

$ws.Range("B29").Value = 6977
$ws.Range("C29").Value = 43933
$ws.Range("D29").Value = 0.84583333333333333
$ws.Range("E29").Value = 0.85069444444444453
$ws.Range("G29").Value = "Ran functional simulations for ShiftUnit.vhd. Shared results of wrong results with team members. DONE"

$ws.Range("B30").Value = 6977
$ws.Range("C30").Value = 43933
$ws.Range("D30").Value = 0.85069444444444453
$ws.Range("E30").Value = 0.85763888888888884
$ws.Range("G30").Value = "Ran functional simulations for ExecUnit.vhd. Shared results of wrong results with team members. DONE"

$ws.Range("B31").Value = 6977
$ws.Range("C31").Value = 43933
$ws.Range("D31").Value = 0.85763888888888884
$ws.Range("E31").Value = 0.86111111111111116
$ws.Range("G31").Value = "Committed project files and code to Github. Taking a break for supper"

# Update the selected cell shown on the sheet
$ws.Activate()
$ws.Range("E20").Select()
